$d = $word.ActiveDocument

# Paragraph 1 = letterhead/address block
# Paragraph 2 = date
# Paragraph 3 = salutation
# Paragraph 4 = first body paragraph (mentions journal section + journal name)
# Paragraph 5 = second body paragraph (mentions journal section + journal name)
#
# NOTE: Range.Find.Execute(..., Replace:=wdReplaceAll) collapses the Range
# to the last replaced text, so we must re-fetch Paragraphs.Item(N).Range
# fresh before every Find call instead of reusing a cached Range variable -
# otherwise the next Find on that "range" only covers the previous match
# and silently falls through to match elsewhere in the document.

# --- Address block (letterhead) ---
$d.Paragraphs.Item(1).Range.Find.Execute("Howard I. Browman, Ph.D.", $true, $false, $false, $false, $false, $true, 1, $false, "Editorial office", 2) | Out-Null
$d.Paragraphs.Item(1).Range.Find.Execute("Editor-in-Chief", $true, $false, $false, $false, $false, $true, 1, $false, "Canadian Journal of Fisheries and Aquatic Sciences", 2) | Out-Null
$d.Paragraphs.Item(1).Range.Find.Execute("ICES Journal of Marine Science", $true, $false, $false, $false, $false, $true, 1, $false, "Canadian Science Publishing (NRC Research Press)", 2) | Out-Null
$d.Paragraphs.Item(1).Range.Find.Execute("Institute of Marine Research", $true, $false, $false, $false, $false, $true, 1, $false, "65 Auriga Drive, Suite 203", 2) | Out-Null
$d.Paragraphs.Item(1).Range.Find.Execute("Marine Ecosystem Acoustics Group", $true, $false, $false, $false, $false, $true, 1, $false, "Ottawa, ON K2E 7W6", 2) | Out-Null

# Remove the "Austevoll Research Station" / "5392 Storebø" lines (and their
# line breaks) entirely - vertical-tab (Chr 11) is how Word represents the
# <w:br w:type="textWrapping"/> manual line break inside Range.Text.
$vt = [char]11
$removeSpan = "Austevoll Research Station" + $vt + "5392 Storebø" + $vt
$d.Paragraphs.Item(1).Range.Find.Execute($removeSpan, $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

$d.Paragraphs.Item(1).Range.Find.Execute("Norway", $true, $false, $false, $false, $false, $true, 1, $false, "Canada", 2) | Out-Null

# --- Date ---
$d.Paragraphs.Item(2).Range.Find.Execute("2019-10-23", $true, $false, $false, $false, $false, $true, 1, $false, "2019-11-30", 2) | Out-Null

# --- Salutation ---
$d.Paragraphs.Item(3).Range.Find.Execute("Dr. Browman,", $true, $false, $false, $false, $false, $true, 1, $false, "Dear Editors,", 2) | Out-Null

# --- Body text: journal section type and journal name (2 occurrences, one per paragraph) ---
$d.Paragraphs.Item(4).Range.Find.Execute("Quo Vadimus", $true, $false, $false, $false, $false, $true, 1, $false, "Perspectives", 2) | Out-Null
$d.Paragraphs.Item(4).Range.Find.Execute("ICES Journal of Marine Science", $true, $false, $false, $false, $false, $true, 1, $false, "Canadian Journal of Fisheries and Aquatic Sciences", 2) | Out-Null

$d.Paragraphs.Item(5).Range.Find.Execute("Quo Vadimus", $true, $false, $false, $false, $false, $true, 1, $false, "Perspectives", 2) | Out-Null
$d.Paragraphs.Item(5).Range.Find.Execute("ICES Journal of Marine Science", $true, $false, $false, $false, $false, $true, 1, $false, "Canadian Journal of Fisheries and Aquatic Sciences", 2) | Out-Null
